$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row text (A1 stays "state_name"; B1 "counties" -> "county_names") ---
$ws.Range("A1").Value2 = "state_name"
$ws.Range("B1").Value2 = "county_names"

# --- Update data rows: append County/city designations to each place name ---
$ws.Range("A2").Value2 = "Maryland"
$ws.Range("B2").Value2 = "Allegany County,Anne Arundel County,Baltimore County,Baltimore city,Calvert County,Caroline County,Carroll County,Cecil County,Charles County,Dorchester County,Frederick County,Garrett County,Harford County,Howard County,Kent County,Montgomery County,Prince George's County,Queen Anne's County,Somerset County,St. Mary's County,Talbot County,Washington County,Wicomico County,Worcester County"

$ws.Range("A3").Value2 = "Virginia"
$ws.Range("B3").Value2 = "Accomack County,Albemarle County,Alexandria city,Alleghany County,Amelia County,Amherst County,Appomattox County,Arlington County,Augusta County,Bath County,Bedford County,Bland County,Botetourt County,Bristol city,Brunswick County,Buchanan County,Buckingham County,Buena Vista city,Campbell County,Caroline County,Carroll County,Charles City County,Charlotte County,Charlottesville city,Chesapeake city,Chesterfield County,Clarke County,Colonial Heights city,Covington city,Craig County,Culpeper County,Cumberland County,Danville city,Dickenson County,Dinwiddie County,Emporia city,Essex County,Fairfax County,Fairfax city,Falls Church city,Fauquier County,Floyd County,Fluvanna County,Franklin County,Franklin city,Frederick County,Fredericksburg city,Galax city,Giles County,Gloucester County,Goochland County,Grayson County,Greene County,Greensville County,Halifax County,Hampton city,Hanover County,Harrisonburg city,Henrico County,Henry County,Highland County,Hopewell city,Isle of Wight County,James City County,King George County,King William County,King and Queen County,Lancaster County,Lee County,Lexington city,Loudoun County,Louisa County,Lunenburg County,Lynchburg city,Madison County,Manassas city,Manassas Park city,Martinsville city,Mathews County,Mecklenburg County,Middlesex County,Montgomery County,Nelson County,New Kent County,Newport News city,Norfolk city,Northampton County,Northumberland County,Norton city,Nottoway County,Orange County,Page County,Patrick County,Petersburg city,Pittsylvania County,Poquoson city,Portsmouth city,Powhatan County,Prince Edward County,Prince George County,Prince William County,Pulaski County,Radford city,Rappahannock County,Richmond County,Richmond city,Roanoke County,Roanoke city,Rockbridge County,Rockingham County,Russell County,Salem city,Scott County,Shenandoah County,Smyth County,Southampton County,Spotsylvania County,Stafford County,Staunton city,Suffolk city,Surry County,Sussex County,Tazewell County,Virginia Beach city,Warren County,Washington County,Waynesboro city,Westmoreland County,Williamsburg city,Winchester city,Wise County,Wythe County,York County"

# --- Remove wrap-text formatting from column B cells (header + data) ---
$ws.Range("B1").WrapText = $false
$ws.Range("B2").Style = "Normal"
$ws.Range("B3").Style = "Normal"

# --- Re-fit rows/columns now that wrapping is gone, restoring default sizing ---
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()

# --- Clear the selection anchor back to the top-left cell ---
[void]$ws.Range("A1").Select()
